$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$cell = $ws.Cells.Item(2, 4)
$cell.NumberFormat = "@"
$cell.Value = "70.827.17"
$cell.Style = "Normal"
$ws.Cells.Item(2, 5).Value = "  -2.09%  "

# Row 3
$cell = $ws.Cells.Item(3, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.553.96"
$cell.Style = "Normal"
$ws.Cells.Item(3, 5).Value = "  -5.82%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  +0.12%  "

# Row 5
$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = "578.11"
$cell.Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  -3.78%  "

# Row 6
$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = "@"
$cell.Value = "169.92"
$cell.Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  -3.68%  "

# Row 7
$ws.Cells.Item(7, 5).Value = "  +0.08%  "

# Row 8
$cell = $ws.Cells.Item(8, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.511"
$cell.Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  -2.74%  "

# Row 9
$ws.Cells.Item(9, 5).Value = "  -1.77%  "

# Row 10
$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.558.29"
$cell.Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  -5.70%  "

# Row 11
$cell = $ws.Cells.Item(11, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.168"
$cell.Style = "Normal"

# Row 12
$cell = $ws.Cells.Item(12, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.350"
$cell.Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  -1.49%  "

# Row 13
$cell = $ws.Cells.Item(13, 4)
$cell.NumberFormat = "@"
$cell.Value = "4.84"
$cell.Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  -3.70%  "

# Row 14
$ws.Cells.Item(14, 2).Value = "WrappedliquidstakedEther2.0"
$ws.Cells.Item(14, 3).Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$cell = $ws.Cells.Item(14, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.054.89"
$cell.Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  -4.75%  "

# Row 15
$ws.Cells.Item(15, 2).Value = "ShibaInu"
$ws.Cells.Item(15, 3).Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$cell = $ws.Cells.Item(15, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0000185"
$cell.Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  -0.67%  "

# Row 16
$cell = $ws.Cells.Item(16, 4)
$cell.NumberFormat = "@"
$cell.Value = "70.725.29"
$cell.Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  -2.02%  "

# Row 17
$cell = $ws.Cells.Item(17, 4)
$cell.NumberFormat = "@"
$cell.Value = "25.18"
$cell.Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  -4.60%  "

# Row 18
$cell = $ws.Cells.Item(18, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.568.83"
$cell.Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  -5.38%  "

# Row 19
$cell = $ws.Cells.Item(19, 4)
$cell.NumberFormat = "@"
$cell.Value = "11.85"
$cell.Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  -3.87%  "

# Row 20
$cell = $ws.Cells.Item(20, 4)
$cell.NumberFormat = "@"
$cell.Value = "7.68"
$cell.Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  -5.84%  "

# Row 21
$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = "@"
$cell.Value = "364.17"
$cell.Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  -2.69%  "

# Row 22
$cell = $ws.Cells.Item(22, 4)
$cell.NumberFormat = "@"
$cell.Value = "4.00"
$cell.Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  -4.29%  "

# Row 23
$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.99"
$cell.Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  -2.01%  "

# Row 24
$cell = $ws.Cells.Item(24, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.999"
$cell.Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  -0.11%  "

# Row 25
$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = "@"
$cell.Value = "70.33"
$cell.Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  -2.99%  "

# Row 26
$cell = $ws.Cells.Item(26, 4)
$cell.NumberFormat = "@"
$cell.Value = "4.13"
$cell.Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  -5.92%  "

# Row 27
$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = "@"
$cell.Value = "9.26"
$cell.Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  -6.01%  "

# Row 28
$cell = $ws.Cells.Item(28, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.701.02"
$cell.Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  -5.35%  "

# Row 29
$cell = $ws.Cells.Item(29, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.999"
$cell.Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  +0.06%  "

# Row 30
$cell = $ws.Cells.Item(30, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0₃0927"
$cell.Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  -7.14%  "

# Row 31
$cell = $ws.Cells.Item(31, 4)
$cell.NumberFormat = "@"
$cell.Value = "7.78"
$cell.Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  -4.50%  "

# Row 32
$cell = $ws.Cells.Item(32, 4)
$cell.NumberFormat = "@"
$cell.Value = "484.65"
$cell.Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  -4.85%  "

# Row 33
$cell = $ws.Cells.Item(33, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.30"
$cell.Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  -1.56%  "

# Row 34
$cell = $ws.Cells.Item(34, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.77"
$cell.Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  -3.10%  "

# Row 35
$ws.Cells.Item(35, 5).Value = "  +0.21%  "

# Row 36
$cell = $ws.Cells.Item(36, 4)
$cell.NumberFormat = "@"
$cell.Value = "157.27"
$cell.Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  -4.17%  "

# Row 37
$cell = $ws.Cells.Item(37, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.114"
$cell.Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  +4.55%  "

# Row 38
$cell = $ws.Cells.Item(38, 4)
$cell.NumberFormat = "@"
$cell.Value = "18.80"
$cell.Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  -4.61%  "

# Row 39
$cell = $ws.Cells.Item(39, 4)
$cell.NumberFormat = "@"
$cell.Value = "18.85"
$cell.Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  -1.41%  "

# Row 40
$cell = $ws.Cells.Item(40, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.33"
$cell.Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  -4.99%  "

# Row 41
$ws.Cells.Item(41, 5).Value = "  -0.08%  "

# Row 42
$cell = $ws.Cells.Item(42, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.69"
$cell.Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  -6.59%  "

# Row 43
$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.47"
$cell.Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  -4.15%  "

# Row 44
$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = "@"
$cell.Value = "4.76"
$cell.Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  -6.32%  "

# Row 45
$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.320"
$cell.Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  -4.91%  "

# Row 46
$cell = $ws.Cells.Item(46, 4)
$cell.NumberFormat = "@"
$cell.Value = "38.61"
$cell.Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  -2.17%  "

# Row 47
$cell = $ws.Cells.Item(47, 4)
$cell.NumberFormat = "@"
$cell.Value = "146.64"
$cell.Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  -6.45%  "

# Row 48
$cell = $ws.Cells.Item(48, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.57"
$cell.Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  -5.10%  "

# Row 49
$cell = $ws.Cells.Item(49, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.530"
$cell.Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  -6.47%  "

# Row 50
$cell = $ws.Cells.Item(50, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.63"
$cell.Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  -8.46%  "

# Row 51
$cell = $ws.Cells.Item(51, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.594"
$cell.Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  -2.45%  "
